{"js": "// Update the worksheet date and every three-digit \u00f7 one-digit division\n// answer cell to the newly generated problems/results. Each old string\n// below is unique within the document body, so a search()+insertText\n// replace per pair unambiguously targets the right run.\nconst replacements = [\n  [\"2025-03-17 Monday\", \"2025-03-18 Tuesday\"],\n  [\"447\u00f74=111, 3\", \"535\u00f76=89, 1\"],\n  [\"990\u00f76=165, 0\", \"863\u00f76=143, 5\"],\n  [\"236\u00f78=29, 4\", \"600\u00f78=75, 0\"],\n  [\"493\u00f76=82, 1\", \"490\u00f77=70, 0\"],\n  [\"538\u00f75=107, 3\", \"228\u00f74=57, 0\"],\n  [\"249\u00f79=27, 6\", \"908\u00f72=454, 0\"],\n  [\"373\u00f73=124, 1\", \"244\u00f73=81, 1\"],\n  [\"491\u00f77=70, 1\", \"913\u00f76=152, 1\"],\n  [\"312\u00f73=104, 0\", \"600\u00f78=75, 0\"],\n  [\"467\u00f72=233, 1\", \"693\u00f77=99, 0\"],\n  [\"172\u00f78=21, 4\", \"973\u00f73=324, 1\"],\n  [\"999\u00f72=499, 1\", \"956\u00f75=191, 1\"],\n  [\"937\u00f77=133, 6\", \"235\u00f75=47, 0\"],\n  [\"918\u00f79=102, 0\", \"247\u00f78=30, 7\"],\n  [\"912\u00f73=304, 0\", \"888\u00f72=444, 0\"],\n  [\"288\u00f75=57, 3\", \"370\u00f77=52, 6\"],\n  [\"605\u00f77=86, 3\", \"957\u00f73=319, 0\"],\n  [\"688\u00f73=229, 1\", \"849\u00f72=424, 1\"],\n  [\"298\u00f75=59, 3\", \"139\u00f74=34, 3\"],\n  [\"784\u00f75=156, 4\", \"587\u00f76=97, 5\"],\n  [\"864\u00f74=216, 0\", \"289\u00f78=36, 1\"],\n  [\"795\u00f76=132, 3\", \"323\u00f74=80, 3\"],\n  [\"774\u00f75=154, 4\", \"432\u00f72=216, 0\"],\n  [\"891\u00f72=445, 1\", \"911\u00f79=101, 2\"],\n  [\"938\u00f75=187, 3\", \"159\u00f78=19, 7\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each three-digit/one-digit division answer\n# with the newly generated values. Each old string is unique in the\n# document, so Find/Replace (ReplaceAll) per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-03-17 Monday\", \"2025-03-18 Tuesday\"),\n    @(\"447\u00f74=111, 3\", \"535\u00f76=89, 1\"),\n    @(\"990\u00f76=165, 0\", \"863\u00f76=143, 5\"),\n    @(\"236\u00f78=29, 4\", \"600\u00f78=75, 0\"),\n    @(\"493\u00f76=82, 1\", \"490\u00f77=70, 0\"),\n    @(\"538\u00f75=107, 3\", \"228\u00f74=57, 0\"),\n    @(\"249\u00f79=27, 6\", \"908\u00f72=454, 0\"),\n    @(\"373\u00f73=124, 1\", \"244\u00f73=81, 1\"),\n    @(\"491\u00f77=70, 1\", \"913\u00f76=152, 1\"),\n    @(\"312\u00f73=104, 0\", \"600\u00f78=75, 0\"),\n    @(\"467\u00f72=233, 1\", \"693\u00f77=99, 0\"),\n    @(\"172\u00f78=21, 4\", \"973\u00f73=324, 1\"),\n    @(\"999\u00f72=499, 1\", \"956\u00f75=191, 1\"),\n    @(\"937\u00f77=133, 6\", \"235\u00f75=47, 0\"),\n    @(\"918\u00f79=102, 0\", \"247\u00f78=30, 7\"),\n    @(\"912\u00f73=304, 0\", \"888\u00f72=444, 0\"),\n    @(\"288\u00f75=57, 3\", \"370\u00f77=52, 6\"),\n    @(\"605\u00f77=86, 3\", \"957\u00f73=319, 0\"),\n    @(\"688\u00f73=229, 1\", \"849\u00f72=424, 1\"),\n    @(\"298\u00f75=59, 3\", \"139\u00f74=34, 3\"),\n    @(\"784\u00f75=156, 4\", \"587\u00f76=97, 5\"),\n    @(\"864\u00f74=216, 0\", \"289\u00f78=36, 1\"),\n    @(\"795\u00f76=132, 3\", \"323\u00f74=80, 3\"),\n    @(\"774\u00f75=154, 4\", \"432\u00f72=216, 0\"),\n    @(\"891\u00f72=445, 1\", \"911\u00f79=101, 2\"),\n    @(\"938\u00f75=187, 3\", \"159\u00f78=19, 7\")\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $rng = $d.Content\n    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n"}
